$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.946.28"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.513.08"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "323.59"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "109.85"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.33%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "41.00"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +7.10%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  +0.65%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "18.75"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.60%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.27"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "2.906.56"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "2.509.86"
$ws.Range("E16").Value = "  +1.02%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.859"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "47.847.42"
$ws.Range("E18").Value = "  +1.18%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.36"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +4.11%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.67"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.80"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +15.50%  "
$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("E22").Value = "  +0.79%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.93"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "248.21"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  +0.11%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "26.01"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.06"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +3.46%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "35.16"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "49.78"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "20.19"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.38"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0791"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +0.22%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "4.70"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  +0.20%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "22.81"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +7.87%  "
$ws.Range("E42").Value = "  -0.94%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "119.65"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.85%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0299"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "2.005.24"
$ws.Range("E45").Value = "  +2.08%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("E49").Value = "  -0.70%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "5.22"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "57.22"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +4.13%  "
